# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and
# apply GaNParagraph / GaNLinks to the relevant runs, per the commit:
# "Add styles to the new paragraphs".

$d = $word.ActiveDocument

# --- Define the new character styles -------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNParagraph to every "Du deltar..." run ------------------------

$paragraphText = "Du deltar i en världsomspännande kampanj för att observera och rapportera de svagaste synliga stjärnorna, som ett mått på ljusföroreningarna på orten. Genom att hitta och observera Orion konstellation på natthimlen kan folk i hela världen lära sig hur belysningen i våra samhällen och omgivningar bidrar till ljusföroreningar. Era bidrag till online-databasen hjälper till att dokumentera den synliga natthimlens över hela världen."

$rng = $d.Content
$rng.Start = 0
while ($rng.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# --- Apply GaNLinks to the amper.ped.muni.cz link run ----------------------

$linkText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng2 = $d.Content
$rng2.Start = 0
while ($rng2.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNLinks"
    $rng2.Collapse(0)
    $rng2.End = $d.Content.End
}

Write-Output "done"
